$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (row 6) following the existing pattern
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = 45674

$ws.Range("B6").Value = "1.25 hours"
$ws.Range("C6").Value = "cleaning data"
$ws.Range("D6").Value = "N/A"

# Update selection to match post-edit state (next empty row)
$ws.Range("B7").Select()
